$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nerfs to Merchant Class Specialties
# Update agi_mod (column L) values for rows 3-10
$ws.Range("L3").Value = 0.1
$ws.Range("L4").Value = 0.2
$ws.Range("L5").Value = 0.29
$ws.Range("L6").Value = 0.38
$ws.Range("L7").Value = 0.47
$ws.Range("L8").Value = 0.57
$ws.Range("L9").Value = 0.66
$ws.Range("L10").Value = 0.75

# Remove the damage column (W) values for rows 2-10
$ws.Range("W2:W10").ClearContents()
